# Update the Monte Carlo simulation output values on Sheet1 (B2:K15).
# The new values replace the previous simulation run's results.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(
    @(0.3311781003110986,0.331810467701082,0.3304907755749598,0.330722976418693,0.329926366626636,0.3319838159440212,0.3317016248101555,0.3304946095861987,0.3301633671745965,0.330372746048705),
    @(37.52318014156442,37.48865765372842,37.62588580448091,37.60321027042146,37.67879401029531,37.44226377268393,37.4900940109327,37.63276016822685,37.68932425269262,37.66121525110508),
    @(709.2371634802125,720.5964926918749,720.1587211052852,721.3219288212722,699.8105848884077,700.4867193667177,729.4380107881594,720.862732101896,695.4736102889991,730.9137665333081),
    @(62.01472133332415,64.04126996005594,65.50414802624466,62.48059010844163,62.60377756747553,61.47926184595418,63.28363141773372,64.20520981716984,61.90707249802217,62.34659459336421),
    @(19885.22541514538,20969.84964751052,21645.88385862309,20268.5015865637,20031.17119339164,19545.95270952017,20724.33714801118,21049.64147953193,19669.13365028307,20308.97274119411),
    @(552.8613174278472,743.5671630108276,1169.931550073584,736.7451579240563,649.2934843495186,652.0519624888124,1276.930316904113,989.86775290715,421.3976863541116,961.5712066029268),
    @(-540.0985721806927,-285.354658246741,236.9922254774391,-425.2078256488876,-182.304109885095,-530.1519616572132,-502.7959479337884,-2171.395772588136,-84.37261637634892,-210.4277362492401),
    @(1801.837703558495,1834.512773557341,1844.959466870111,1805.166853164672,1811.312930562066,1791.037257972854,1815.539153478961,1836.42777183439,1786.961449334356,1812.61941764426),
    @(4420.459186295448,8739.098584772299,8307.633925786151,8671.429491083627,17225.35086019996,8553.245754893975,15030.47197700065,17257.1603460054,16822.58110952224,5088.640144949332),
    @(5870.87660356945,7675.205870600206,9110.348961400712,13100.31674337607,6307.160926262772,5358.123224318038,13915.34514941047,6533.412475188299,12032.00273418819,3494.037436894518),
    @(-4.027782221470205,-4.013623046134256,-4.041099256203577,-4.041654684971183,-4.014833822682297,-4.025666737424493,-4.066179095153172,-4.00457905950386,-4.081346869426931,-4.028905620320212),
    @(-1.381446904341689,-1.729028545916783,-1.655969407291891,-1.714239868074755,-1.252697993676119,-1.581199090012658,-1.832524024385572,-1.144894795230538,-1.646633832369272,-1.452318890346519),
    @(-1.327623023796891,-1.675286892609603,-1.602767056219131,-1.660827960934341,-1.199621683980986,-1.526941373389958,-1.778634520752095,-1.091779081561023,-1.593561433023413,-1.39941489554721),
    @(2.09286643248244,1.737566141289081,1.777843532012253,1.733417496748528,2.175033380208535,1.919411630355777,1.644069335336403,2.283189166746962,1.778633136110408,1.962016874373689)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $rowIndex = $i + 2   # data starts at row 2
    $rowVals = $values[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $colIndex = $j + 2   # data starts at column B
        $ws.Cells.Item($rowIndex, $colIndex).Value = $rowVals[$j]
    }
}
